$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New email / password rows appended below the existing data (rows 3-10).
$rows = @(
    @{ Row = 3;  Email = "ananthram1113@gmail.com";     Password = "Ananth@4a7" },
    @{ Row = 4;  Email = "adithyaram1113@gmail.com";    Password = "Adithya@113" },
    @{ Row = 5;  Email = "srikanthreddy459@gmail.com";  Password = "Srikanth@12" },
    @{ Row = 6;  Email = "Yeswanth231@gmail.com";       Password = "Yeswanth@451" },
    @{ Row = 7;  Email = "Hema.selenium342@gmail.com";  Password = "Hema@342" },
    @{ Row = 8;  Email = "naveensai649@gmail.com";      Password = "Naveen@649" },
    @{ Row = 9;  Email = "ramprasad.selenium@gmail.com";Password = "Prasad@123" },
    @{ Row = 10; Email = "majjiteja.98@gmail.com";      Password = "Teja@98" }
)

foreach ($r in $rows) {
    $emailCell = $ws.Cells.Item($r.Row, 1)
    $passwordCell = $ws.Cells.Item($r.Row, 2)

    $emailCell.Value = $r.Email
    $passwordCell.Value = $r.Password

    [void]$ws.Hyperlinks.Add($emailCell, "mailto:" + $r.Email)
    [void]$ws.Hyperlinks.Add($passwordCell, "mailto:" + $r.Password)
}

# The Hyperlinks.Add calls stamp their own cell style; bring every data row
# back to the shared "Hyperlink" cell style (same as the existing row 2).
$ws.Range("A3:B10").Style = "Hyperlink"

# Match the column widths / layout nudge that shipped with the new rows.
$ws.Columns.Item(1).ColumnWidth = 33.5
$ws.Columns.Item(2).ColumnWidth = 15.5

# Move the active selection past the last populated row, like the source file.
[void]$ws.Range("B11").Select()
